$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1: "Age" -> "City"
$ws.Range("C1").Value = "City"

# Update data rows: column C values change from numeric ages to city name strings
$ws.Range("C2").Value = "Bharatpur"
$ws.Range("C3").Value = "Bangalore"

# Update selection to C3 as in the edited file
$ws.Range("C3").Select()
